# Add a new "2022" column (P) to the human-trafficking-victims table,
# mirroring the formatting of the existing "2021" column (O), then fill
# in the 2022 figures, and finally move the selection to match the
# author's saved state (O21:O22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone formatting (styles/number formats/borders) from column O (2021)
# into column P for the header + all data rows, including the bottom
# "totals" row.
$ws.Range("O4:O14").Copy($ws.Range("P4")) | Out-Null
$excel.CutCopyMode = 0

# Header: new year
$ws.Range("P4").Value = 2022

# Data rows for 2022
$ws.Range("P5").Value = 1
$ws.Range("P6").Value = "-"
$ws.Range("P7").Value = "-"
$ws.Range("P8").Value = "-"
$ws.Range("P9").Value = "-"
$ws.Range("P10").Value = "-"
$ws.Range("P11").Value = "-"
$ws.Range("P12").Value = 1
$ws.Range("P13").Value = "-"
$ws.Range("P14").Value = "-"

# Move the active selection to match the final saved workbook state.
$ws.Range("O21:O22").Select() | Out-Null
